$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ExistingRunConfig_Data")

$ws.Range("A16").Value = "Level 16"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "NO"

$ws.Range("A17").Value = "Level 17"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "NO"

$ws.Range("A18").Value = "Level 18"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "NO"

$ws.Range("A19").Value = "Level 19"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = "NO"

$ws.Range("A20").Value = "Level 20"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = "NO"

$ws.Range("A21").Value = "Level 21"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "NO"

$ws.Range("A22").Value = "Level 22"
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = "NO"

$ws.Range("A23").Value = "Level 23"
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 24
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = "NO"

$ws.Range("A24").Value = "Level 24"
$ws.Range("B24").Value = 2
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = "NO"

$ws.Range("A25").Value = "Level 25"
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = "NO"

$ws.Range("A26").Value = "Level 26"
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 4
$ws.Range("F26").Value = "NO"

$ws.Range("A27").Value = "Level 27"
$ws.Range("B27").Value = 1
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 1
$ws.Range("F27").Value = "NO"

$ws.Range("A28").Value = "Level 28"
$ws.Range("B28").Value = 99
$ws.Range("C28").Value = 19
$ws.Range("D28").Value = 19
$ws.Range("E28").Value = 19
$ws.Range("F28").Value = "NO"

$ws.Range("A29").Value = "Level 28"
$ws.Range("B29").Value = 99
$ws.Range("C29").Value = 19
$ws.Range("D29").Value = 19
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = "NO"

$ws.Range("A30").Value = "Level 29"
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 2
$ws.Range("E30").Value = 2
$ws.Range("F30").Value = "NO"

$ws.Range("A31").Value = "Level 30"
$ws.Range("B31").Value = 3
$ws.Range("C31").Value = 2
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = "NO"

$ws.Range("A32").Value = "Level 31"
$ws.Range("B32").Value = 11
$ws.Range("C32").Value = 11
$ws.Range("D32").Value = 11
$ws.Range("E32").Value = 11
$ws.Range("F32").Value = "NO"

$ws.Range("A33").Value = "Level 32"
$ws.Range("B33").Value = 11
$ws.Range("C33").Value = 11
$ws.Range("D33").Value = 11
$ws.Range("E33").Value = 11
$ws.Range("F33").Value = "NO"

$ws.Range("A34").Value = "Level 33"
$ws.Range("B34").Value = 1
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 1
$ws.Range("F34").Value = "NO"

$ws.Range("A35").Value = "Level 34"
$ws.Range("B35").Value = 3
$ws.Range("C35").Value = 7
$ws.Range("D35").Value = 5
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = "NO"

$ws.Range("A36").Value = "Level 35"
$ws.Range("B36").Value = 1
$ws.Range("C36").Value = 8
$ws.Range("D36").Value = 5
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = "NO"

$ws.Range("A37").Value = "Level 36"
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = "NO"

$ws.Range("A38").Value = "Level 37"
$ws.Range("B38").Value = 2
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = "NO"
